# Started On sorting of ColumnHeader functionality
# Adds a new "verifyContactHeaders" worksheet (becomes the active/selected
# tab) that lists the contact table's column headers.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet ("NewEvent") so it lands
# at the end of the tab strip, then name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "verifyContactHeaders"

# Header label (bold) followed by a blank row, then each column heading.
$newSheet.Range("A1").Value = "ContactHeader"
$newSheet.Range("A1").Font.Bold = $true

$newSheet.Range("A3").Value = "Name"
$newSheet.Range("A4").Value = "Address"
$newSheet.Range("A5").Value = "Category"
$newSheet.Range("A6").Value = "Status"
$newSheet.Range("A7").Value = "Phone"
$newSheet.Range("A8").Value = "Email"
$newSheet.Range("A9").Value = "Options"

# Match the source column width as closely as possible and set the
# selection to where the author last left the cursor on this sheet.
$newSheet.Columns.Item(1).ColumnWidth = 12.330729166666666
[void]$newSheet.Range("H35").Select()
